# Rename sheets and delete the obsolete "Desarquivamentos Pendentes" sheet.

$wb = $excel.ActiveWorkbook

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet entirely
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
